# Update "Pais" sheet: refresh country stats and fix rank ordering
# (countries swap position in the shared-string table as case counts change)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 'Estados Unidos'
$ws.Range("B4").Value = 1352313
$ws.Range("C4").Value = 5004
$ws.Range("D4").Value = 238848
$ws.Range("E4").Value = 1033279
$ws.Range("F4").Value = 16816
$ws.Range("G4").Value = 149
$ws.Range("H4").Value = 80186

$ws.Range("A6").Value = 'Reino Unido'
$ws.Range("B6").Value = 219183
$ws.Range("C6").Value = 3923
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 186984
$ws.Range("F6").Value = 1559
$ws.Range("G6").Value = 268
$ws.Range("H6").Value = 31855

$ws.Range("A7").Value = 'Italia'
$ws.Range("B7").Value = 219070
$ws.Range("C7").Value = 802
$ws.Range("D7").Value = 105186
$ws.Range("E7").Value = 83324
$ws.Range("F7").Value = 1027
$ws.Range("G7").Value = 165
$ws.Range("H7").Value = 30560

$ws.Range("A10").Value = 'Alemania'
$ws.Range("B10").Value = 171639
$ws.Range("C10").Value = 315
$ws.Range("D10").Value = 144400
$ws.Range("E10").Value = 19690
$ws.Range("F10").Value = 1650
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 7549

$ws.Range("A12").Value = 'Turquia'
$ws.Range("B12").Value = 138657
$ws.Range("C12").Value = 1542
$ws.Range("D12").Value = 92691
$ws.Range("E12").Value = 42180
$ws.Range("F12").Value = 1154
$ws.Range("G12").Value = 47
$ws.Range("H12").Value = 3786

$ws.Range("A22").Value = 'Pakistan'
$ws.Range("B22").Value = 30334
$ws.Range("C22").Value = 1598
$ws.Range("D22").Value = 8023
$ws.Range("E22").Value = 21652
$ws.Range("F22").Value = 111
$ws.Range("G22").Value = 23
$ws.Range("H22").Value = 659

$ws.Range("A23").Value = 'Suiza'
$ws.Range("B23").Value = 30305
$ws.Range("C23").Value = 54
$ws.Range("D23").Value = 26400
$ws.Range("E23").Value = 2075
$ws.Range("F23").Value = 101
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 1830

$ws.Range("A28").Value = 'Singapur'
$ws.Range("B28").Value = 23336
$ws.Range("C28").Value = 876
$ws.Range("D28").Value = 2721
$ws.Range("E28").Value = 20595
$ws.Range("F28").Value = 23
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 20

$ws.Range("A34").Value = 'Polonia'
$ws.Range("B34").Value = 15996
$ws.Range("C34").Value = 345
$ws.Range("D34").Value = 5698
$ws.Range("E34").Value = 9498
$ws.Range("F34").Value = 160
$ws.Range("G34").Value = 15
$ws.Range("H34").Value = 800

$ws.Range("A35").Value = 'Austria'
$ws.Range("B35").Value = 15871
$ws.Range("C35").Value = 38
$ws.Range("D35").Value = 13991
$ws.Range("E35").Value = 1262
$ws.Range("F35").Value = 72
$ws.Range("G35").Value = 3
$ws.Range("H35").Value = 618

$ws.Range("A37").Value = 'Rumania'
$ws.Range("B37").Value = 15362
$ws.Range("C37").Value = 231
$ws.Range("D37").Value = 7051
$ws.Range("E37").Value = 7350
$ws.Range("F37").Value = 242
$ws.Range("G37").Value = 22
$ws.Range("H37").Value = 961

$ws.Range("A45").Value = 'Republica Dominicana'
$ws.Range("B45").Value = 10347
$ws.Range("C45").Value = 465
$ws.Range("D45").Value = 2763
$ws.Range("E45").Value = 7196
$ws.Range("F45").Value = 134
$ws.Range("G45").Value = 3
$ws.Range("H45").Value = 388

$ws.Range("A46").Value = 'Serbia'
$ws.Range("B46").Value = 10114
$ws.Range("C46").Value = 82
$ws.Range("D46").Value = 3006
$ws.Range("E46").Value = 6893
$ws.Range("F46").Value = 43
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = 215

$ws.Range("A48").Value = 'Egipto'
$ws.Range("B48").Value = 9400
$ws.Range("C48").Value = 436
$ws.Range("D48").Value = 2075
$ws.Range("E48").Value = 6800
$ws.Range("F48").Value = 41
$ws.Range("G48").Value = 11
$ws.Range("H48").Value = 525

$ws.Range("A51").Value = 'Chequia'
$ws.Range("B51").Value = 8106
$ws.Range("C51").Value = 11
$ws.Range("D51").Value = 4466
$ws.Range("E51").Value = 3361
$ws.Range("F51").Value = 40
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 279

$ws.Range("A52").Value = 'Noruega'
$ws.Range("B52").Value = 8102
$ws.Range("C52").Value = 3
$ws.Range("D52").Value = 32
$ws.Range("E52").Value = 7851
$ws.Range("F52").Value = 22
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 219

$ws.Range("A60").Value = 'Moldavia'
$ws.Range("B60").Value = 4927
$ws.Range("C60").Value = 60
$ws.Range("D60").Value = 1958
$ws.Range("E60").Value = 2800
$ws.Range("F60").Value = 237
$ws.Range("G60").Value = 8
$ws.Range("H60").Value = 169

$ws.Range("A65").Value = 'Luxemburgo'
$ws.Range("B65").Value = 3886
$ws.Range("C65").Value = 9
$ws.Range("D65").Value = 3586
$ws.Range("E65").Value = 199
$ws.Range("F65").Value = 14
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 101

$ws.Range("A98").Value = 'Consejo Danes para los Refugiados'
$ws.Range("B98").Value = 991
$ws.Range("C98").Value = 54
$ws.Range("D98").Value = 136
$ws.Range("E98").Value = 814
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = 41

$ws.Range("A99").Value = 'Guatemala'
$ws.Range("B99").Value = 967
$ws.Range("C99").Value = 67
$ws.Range("D99").Value = 104
$ws.Range("E99").Value = 839
$ws.Range("F99").Value = 5
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 24

$ws.Range("A100").Value = 'Letonia'
$ws.Range("B100").Value = 939
$ws.Range("C100").Value = 9
$ws.Range("D100").Value = 464
$ws.Range("E100").Value = 457
$ws.Range("F100").Value = 2
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 18

$ws.Range("A104").Value = 'Sri Lanka'
$ws.Range("B104").Value = 856
$ws.Range("C104").Value = 9
$ws.Range("D104").Value = 321
$ws.Range("E104").Value = 526
$ws.Range("F104").Value = 1
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 9

$ws.Range("A111").Value = 'Guinea-Bisau'
$ws.Range("B111").Value = 726
$ws.Range("C111").Value = 85
$ws.Range("D111").Value = 26
$ws.Range("E111").Value = 697
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 3

$ws.Range("A112").Value = 'Crucero'
$ws.Range("B112").Value = 712
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 645
$ws.Range("E112").Value = 54
$ws.Range("F112").Value = 4
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 13

$ws.Range("A113").Value = 'Uruguay'
$ws.Range("B113").Value = 702
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 513
$ws.Range("E113").Value = 171
$ws.Range("F113").Value = 8
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 18

$ws.Range("A114").Value = 'Mali'
$ws.Range("B114").Value = 692
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 298
$ws.Range("E114").Value = 357
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 37

$ws.Range("A115").Value = 'Paraguay'
$ws.Range("B115").Value = 689
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 155
$ws.Range("E115").Value = 524
$ws.Range("F115").Value = 9
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 10

$ws.Range("A116").Value = 'Kenia'
$ws.Range("B116").Value = 672
$ws.Range("C116").Value = 23
$ws.Range("D116").Value = 239
$ws.Range("E116").Value = 401
$ws.Range("F116").Value = 1
$ws.Range("G116").Value = 2
$ws.Range("H116").Value = 32

$ws.Range("A117").Value = 'Gabon'
$ws.Range("B117").Value = 661
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 110
$ws.Range("E117").Value = 543
$ws.Range("F117").Value = 1
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 8

$ws.Range("A132").Value = 'Montenegro'
$ws.Range("B132").Value = 324
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 290
$ws.Range("E132").Value = 25
$ws.Range("F132").Value = 2
$ws.Range("G132").Value = 1
$ws.Range("H132").Value = 9

$ws.Range("A140").Value = 'Cabo Verde'
$ws.Range("B140").Value = 246
$ws.Range("C140").Value = 10
$ws.Range("D140").Value = 56
$ws.Range("E140").Value = 188
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 2

$ws.Range("A141").Value = 'Etiopia'
$ws.Range("B141").Value = 239
$ws.Range("C141").Value = 29
$ws.Range("D141").Value = 99
$ws.Range("E141").Value = 135
$ws.Range("F141").Value = 1
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 5

$ws.Range("A192").Value = 'Nueva Caledonia'
$ws.Range("B192").Value = 18
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 18
$ws.Range("E192").Value = 0
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 0

$ws.Range("A193").Value = 'Belice'
$ws.Range("B193").Value = 18
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 16
$ws.Range("E193").Value = 0
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 2

$ws.Range("A212").Value = 'Islas Virgenes Britanicas'
$ws.Range("B212").Value = 7
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 4
$ws.Range("E212").Value = 2
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 1

$ws.Range("A213").Value = 'Butan'
$ws.Range("B213").Value = 7
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 5
$ws.Range("E213").Value = 2
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

